$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents of C2 and D2 (were "2." and "3.")
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()

# Clear contents of C3 and D3 (were "18°" and "20°")
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()

# Update B4 value, clear C4 and D4
$ws.Range("B4").Value = "09:58:27"
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

# Update B5 value, clear C5 and D5
$ws.Range("B5").Value = "21-09-22"
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
